$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Grand Total row (row 2)
$ws.Range("D2").Value = 46.0
$ws.Range("E2").Value = 16680.78

$ws.Range("A3").Value = "DZCT"
$ws.Range("B3").Value = "MB498897"
$ws.Range("C3").Value = "2025-10-24"
$ws.Range("D3").Value = 1.0
$ws.Range("E3").Value = 711.0

$ws.Range("A4").Value = "DZZ9"
$ws.Range("B4").Value = "MB498511"
$ws.Range("C4").Value = "2025-10-21"
$ws.Range("D4").Value = 1.0
$ws.Range("E4").Value = 382.5

$ws.Range("A5").Value = "DZZ2"
$ws.Range("B5").Value = "MB498257"
$ws.Range("C5").Value = "2025-10-20"
$ws.Range("D5").Value = 1.0
$ws.Range("E5").Value = 334.8

$ws.Range("A6").Value = "DZCT"
$ws.Range("B6").Value = "MB497967"
$ws.Range("C6").Value = "2025-10-19"
$ws.Range("D6").Value = 1.0
$ws.Range("E6").Value = 297.0

$ws.Range("A7").Value = "DZZ2"
$ws.Range("B7").Value = "MB498046"
$ws.Range("C7").Value = "2025-10-19"
$ws.Range("D7").Value = 1.0
$ws.Range("E7").Value = 440.1

$ws.Range("A8").Value = "DZZ2"
$ws.Range("B8").Value = "MB497196"
$ws.Range("C8").Value = "2025-10-15"
$ws.Range("D8").Value = 1.0
$ws.Range("E8").Value = 297.0

$ws.Range("A9").Value = "DZ77"
$ws.Range("B9").Value = "MB497034"
$ws.Range("C9").Value = "2025-10-14"
$ws.Range("D9").Value = 1.0
$ws.Range("E9").Value = 537.3

$ws.Range("A10").Value = "DZ65"
$ws.Range("B10").Value = "MB497171"
$ws.Range("C10").Value = "2025-10-14"
$ws.Range("D10").Value = 1.0
$ws.Range("E10").Value = 359.1

$ws.Range("A11").Value = "DZZ2"
$ws.Range("B11").Value = "MB496903"
$ws.Range("C11").Value = "2025-10-13"
$ws.Range("D11").Value = 1.0
$ws.Range("E11").Value = 297.0

$ws.Range("A12").Value = "DZ77"
$ws.Range("B12").Value = "MB496773"
$ws.Range("C12").Value = "2025-10-12"
$ws.Range("D12").Value = 1.0
$ws.Range("E12").Value = 297.0

$ws.Range("A13").Value = "DZCT"
$ws.Range("B13").Value = "MB496572"
$ws.Range("C13").Value = "2025-10-12"
$ws.Range("D13").Value = 1.0
$ws.Range("E13").Value = 382.5

$ws.Range("A14").Value = "DZCT"
$ws.Range("B14").Value = "MB496514"
$ws.Range("C14").Value = "2025-10-11"
$ws.Range("D14").Value = 1.0
$ws.Range("E14").Value = 297.0

$ws.Range("A15").Value = "DZ77"
$ws.Range("B15").Value = "MB496481"
$ws.Range("C15").Value = "2025-10-11"
$ws.Range("D15").Value = 1.0
$ws.Range("E15").Value = 297.0

$ws.Range("A16").Value = "DZCT"
$ws.Range("B16").Value = "MB496037"
$ws.Range("C16").Value = "2025-10-09"
$ws.Range("D16").Value = 1.0
$ws.Range("E16").Value = 381.78

$ws.Range("A17").Value = "DZZ2"
$ws.Range("B17").Value = "MB496052"
$ws.Range("C17").Value = "2025-10-09"
$ws.Range("D17").Value = 1.0
$ws.Range("E17").Value = 297.0

$ws.Range("A18").Value = "DZCT"
$ws.Range("B18").Value = "MB495517"
$ws.Range("C18").Value = "2025-10-07"
$ws.Range("D18").Value = 1.0
$ws.Range("E18").Value = 486.0

$ws.Range("A19").Value = "DZ65"
$ws.Range("B19").Value = "MB495647"
$ws.Range("C19").Value = "2025-10-07"
$ws.Range("D19").Value = 1.0
$ws.Range("E19").Value = 636.3

$ws.Range("A20").Value = "DZ77"
$ws.Range("B20").Value = "MB495599"
$ws.Range("C20").Value = "2025-10-07"
$ws.Range("D20").Value = 1.0
$ws.Range("E20").Value = 194.4

$ws.Range("A21").Value = "DZ65"
$ws.Range("B21").Value = "MB495453"
$ws.Range("C21").Value = "2025-10-06"
$ws.Range("D21").Value = 1.0
$ws.Range("E21").Value = 269.1

$ws.Range("A22").Value = "DZZ2"
$ws.Range("B22").Value = "MB495050"
$ws.Range("C22").Value = "2025-10-05"
$ws.Range("D22").Value = 2.0
$ws.Range("E22").Value = 703.8

$ws.Range("A23").Value = "DZZ2"
$ws.Range("B23").Value = "MB494491"
$ws.Range("C23").Value = "2025-10-05"
$ws.Range("D23").Value = 1.0
$ws.Range("E23").Value = 351.9

$ws.Range("A24").Value = "DZZ2"
$ws.Range("B24").Value = "MB494936"
$ws.Range("C24").Value = "2025-10-05"
$ws.Range("D24").Value = 1.0
$ws.Range("E24").Value = 495.0

$ws.Range("A25").Value = "DZZ2"
$ws.Range("B25").Value = "MB494715"
$ws.Range("C25").Value = "2025-10-05"
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 267.3

$ws.Range("A26").Value = "DZCT"
$ws.Range("B26").Value = "MB494311"
$ws.Range("C26").Value = "2025-10-04"
$ws.Range("D26").Value = 1.0
$ws.Range("E26").Value = 534.6

$ws.Range("A27").Value = "DZZ2"
$ws.Range("B27").Value = "MB494422"
$ws.Range("C27").Value = "2025-10-04"
$ws.Range("D27").Value = 1.0
$ws.Range("E27").Value = 267.3

$ws.Range("A28").Value = "DZCT"
$ws.Range("B28").Value = "MB493162"
$ws.Range("C28").Value = "2025-10-03"
$ws.Range("D28").Value = 1.0
$ws.Range("E28").Value = 424.8

$ws.Range("A29").Value = "DZ65"
$ws.Range("B29").Value = "MB493466"
$ws.Range("C29").Value = "2025-10-03"
$ws.Range("D29").Value = 1.0
$ws.Range("E29").Value = 246.6

$ws.Range("A30").Value = "DZZ2"
$ws.Range("B30").Value = "MB491923"
$ws.Range("C30").Value = "2025-10-02"
$ws.Range("D30").Value = 1.0
$ws.Range("E30").Value = 267.3

$ws.Range("A31").Value = "DZCT"
$ws.Range("B31").Value = "MB492968"
$ws.Range("C31").Value = "2025-10-02"
$ws.Range("D31").Value = 1.0
$ws.Range("E31").Value = 437.4

$ws.Range("A32").Value = "DZ65"
$ws.Range("B32").Value = "MB492884"
$ws.Range("C32").Value = "2025-10-02"
$ws.Range("D32").Value = 1.0
$ws.Range("E32").Value = 424.8

$ws.Range("A33").Value = "DZCT"
$ws.Range("B33").Value = "MB491721"
$ws.Range("C33").Value = "2025-10-01"
$ws.Range("D33").Value = 1.0
$ws.Range("E33").Value = 267.3

$ws.Range("A34").Value = "DZCT"
$ws.Range("B34").Value = "MB491606"
$ws.Range("C34").Value = "2025-10-01"
$ws.Range("D34").Value = 1.0
$ws.Range("E34").Value = 441.0

$ws.Range("A35").Value = "DZZ2"
$ws.Range("B35").Value = "MB490643"
$ws.Range("C35").Value = "2025-09-30"
$ws.Range("D35").Value = 1.0
$ws.Range("E35").Value = 170.1

$ws.Range("A36").Value = "DZ65"
$ws.Range("B36").Value = "MB490199"
$ws.Range("C36").Value = "2025-09-29"
$ws.Range("D36").Value = 1.0
$ws.Range("E36").Value = 343.8

$ws.Range("A37").Value = "DZZ2"
$ws.Range("B37").Value = "MB489618"
$ws.Range("C37").Value = "2025-09-29"
$ws.Range("D37").Value = 1.0
$ws.Range("E37").Value = 267.3

$ws.Range("A38").Value = "DZZ2"
$ws.Range("B38").Value = "MB489648"
$ws.Range("C38").Value = "2025-09-29"
$ws.Range("D38").Value = 1.0
$ws.Range("E38").Value = 437.4

$ws.Range("A39").Value = "DZCT"
$ws.Range("B39").Value = "MB489754"
$ws.Range("C39").Value = "2025-09-29"
$ws.Range("D39").Value = 1.0
$ws.Range("E39").Value = 522.0

$ws.Range("A40").Value = "DZCT"
$ws.Range("B40").Value = "MB489048"
$ws.Range("C40").Value = "2025-09-28"
$ws.Range("D40").Value = 1.0
$ws.Range("E40").Value = 267.3

$ws.Range("A41").Value = "DZZ2"
$ws.Range("B41").Value = "MB489154"
$ws.Range("C41").Value = "2025-09-28"
$ws.Range("D41").Value = 1.0
$ws.Range("E41").Value = 267.3

$ws.Range("A42").Value = "DZZ2"
$ws.Range("B42").Value = "MB488679"
$ws.Range("C42").Value = "2025-09-28"
$ws.Range("D42").Value = 1.0
$ws.Range("E42").Value = 343.8

$ws.Range("A43").Value = "DZCT"
$ws.Range("B43").Value = "MB487207"
$ws.Range("C43").Value = "2025-09-26"
$ws.Range("D43").Value = 1.0
$ws.Range("E43").Value = 343.8

$ws.Range("A44").Value = "DZZ2"
$ws.Range("B44").Value = "MB487787"
$ws.Range("C44").Value = "2025-09-26"
$ws.Range("D44").Value = 1.0
$ws.Range("E44").Value = 343.8

$ws.Range("A45").Value = "DZZ2"
$ws.Range("B45").Value = "MB487512"
$ws.Range("C45").Value = "2025-09-26"
$ws.Range("D45").Value = 1.0
$ws.Range("E45").Value = 172.8

$ws.Range("A46").Value = "DZCT"
$ws.Range("B46").Value = "MB487316"
$ws.Range("C46").Value = "2025-09-26"
$ws.Range("D46").Value = 1.0
$ws.Range("E46").Value = 611.1

$ws.Range("A47").Value = "DZZ2"
$ws.Range("B47").Value = "MB486945"
$ws.Range("C47").Value = "2025-09-25"
$ws.Range("D47").Value = 1.0
$ws.Range("E47").Value = 267.3

# Clear rows 48 and 49 (no longer have data)
$ws.Range("A48:E49").ClearContents()
